{"js": "// The use-case table's \"Kurzbeschreibung\" (short description) row is empty;\n// fill it in, and the editing cursor's \"_GoBack\" bookmark moves there from the\n// last edited cell (\"Highscore abrufen\"), so that old bookmark must go away.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// --- 1) \"Kurzbeschreibung\" row is the 2nd row (index 1); its 2nd cell\n//        (index 1) currently holds an empty paragraph. ---\nconst descCell = table.getCell(1, 1);\nconst descBody = descCell.body;\ndescBody.load(\"paragraphs/items\");\nawait context.sync();\n\nconst descPara = descBody.paragraphs.items[0];\ndescPara.insertText(\"Der Spieler startet das Spiel\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Drop a \"_GoBack\" bookmark right after the new text (mirrors Word leaving\n// its \"last edit\" bookmark at the most recently typed location).\nconst descEnd = descPara.getRange(Word.RangeLocation.end);\ndescEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) The old \"_GoBack\" bookmark previously sat after \"Highscore abrufen\"\n//        (the last populated cell) and must be removed from there. ---\nconst lastCell = table.getCell(table.rowCount - 1, 1);\nconst lastBody = lastCell.body;\nlastBody.load(\"paragraphs/items\");\nawait context.sync();\n\nconst lastPara = lastBody.paragraphs.items[0];\nlastPara.load(\"text\");\nawait context.sync();\nconst originalText = lastPara.text;\n\n// Collapse-delete the paragraph's content (this removes the bookmark along\n// with the run), then retype the same text so the visible content is\n// unchanged but the bookmark is gone.\nconst lastRange = lastPara.getRange(Word.RangeLocation.end);\nlastRange.delete();\nawait context.sync();\n\nlastPara.insertText(originalText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The use-case table's \"Kurzbeschreibung\" (short description) row is empty;\n# fill it in, and the editor's \"_GoBack\" bookmark moves there from the last\n# edited cell (\"Highscore abrufen\"), so that old bookmark must be removed.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# --- 1) Remove the stale \"_GoBack\" bookmark sitting after \"Highscore\n#        abrufen\" (the table's last populated cell). Deleting through a\n#        manually built Range leaves the bookmark behind, but deleting\n#        through the paragraph's own Range removes it; then retype the\n#        text that was there. ---\n$lastRow = $t.Rows.Count\n$lastCell = $t.Cell($lastRow, 2)\n$lastCell.Range.Paragraphs(1).Range.Delete()\n$t.Cell($lastRow, 2).Range.Text = \"Highscore abrufen\"\n\n# --- 2) Fill in the previously empty \"Kurzbeschreibung\" cell (row 2, col 2). ---\n$descCell = $t.Cell(2, 2)\n$descCell.Range.Text = \"Der Spieler startet das Spiel\"\n\n# --- 3) Re-create the \"_GoBack\" bookmark, collapsed right after the new\n#        text. A zero-length Range placed exactly at a paragraph boundary\n#        is rejected, so temporarily append a throw-away character, anchor\n#        the bookmark just before it (now an interior position, which is\n#        accepted), then delete the throw-away character; the bookmark\n#        stays collapsed exactly where we want it. ---\n$descCell2 = $t.Cell(2, 2)\n$descCell2.Range.InsertAfter(\"X\")\n$pr = $descCell2.Range.Paragraphs(1).Range\n$bmPos = $pr.End - 2\n$bmRange = $d.Range($bmPos, $bmPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n$d.Range($pr.End - 2, $pr.End - 1).Delete()\n"}
